$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update row 4 ("tumor_size") descriptive statistics with newly recomputed values
$ws.Range("D4").Value = 26.05193919343433
$ws.Range("E4").Value = 15.29439420494566
$ws.Range("F4").Value = 0.9497446270443026
$ws.Range("G4").Value = 17
$ws.Range("H4").Value = 22.16989285697366
$ws.Range("K4").Value = 3.043281450982875
$ws.Range("L4").Value = 17.47739465223725
